$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '90.801.46'
$ws.Cells.Item(2, 5).Value = '  +1.57%  '

$ws.Cells.Item(3, 4).Value = '3.162.24'
$ws.Cells.Item(3, 5).Value = '  +4.08%  '

$ws.Cells.Item(5, 4).Value = '215.41'
$ws.Cells.Item(5, 5).Value = '  +2.27%  '

$ws.Cells.Item(6, 4).Value = '627.35'
$ws.Cells.Item(6, 5).Value = '  +2.49%  '

$ws.Cells.Item(7, 5).Value = '  +27.63%  '

$ws.Cells.Item(8, 4).Value = '0.373'
$ws.Cells.Item(8, 5).Value = '  +3.84%  '

$ws.Cells.Item(9, 5).Value = '  -0.02%  '

$ws.Cells.Item(10, 4).Value = '3.159.84'
$ws.Cells.Item(10, 5).Value = '  +4.15%  '

$ws.Cells.Item(11, 4).Value = '0.755'
$ws.Cells.Item(11, 5).Value = '  +12.07%  '

$ws.Cells.Item(12, 5).Value = '  +8.64%  '

$ws.Cells.Item(13, 4).Value = '5.70'
$ws.Cells.Item(13, 5).Value = '  +6.14%  '

$ws.Cells.Item(14, 5).Value = '  +3.03%  '

$ws.Cells.Item(15, 5).Value = '  +8.94%  '

$ws.Cells.Item(16, 4).Value = '90.543.79'
$ws.Cells.Item(16, 5).Value = '  +1.18%  '

$ws.Cells.Item(17, 4).Value = '3.734.66'
$ws.Cells.Item(17, 5).Value = '  +4.01%  '

$ws.Cells.Item(18, 4).Value = '3.141.22'
$ws.Cells.Item(18, 5).Value = '  +3.15%  '

$ws.Cells.Item(19, 5).Value = '  +11.51%  '

$ws.Cells.Item(20, 4).Value = '14.31'
$ws.Cells.Item(20, 5).Value = '  +7.13%  '

$ws.Cells.Item(21, 2).Value = 'PEPE'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(21, 4).Value = '0.0000212'
$ws.Cells.Item(21, 5).Value = '  -1.76%  '

$ws.Cells.Item(22, 2).Value = 'BitcoinCash'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(22, 4).Value = '464.65'
$ws.Cells.Item(22, 5).Value = '  +9.77%  '

$ws.Cells.Item(23, 4).Value = '9.08'
$ws.Cells.Item(23, 5).Value = '  +10.65%  '

$ws.Cells.Item(24, 5).Value = '  +5.69%  '

$ws.Cells.Item(25, 4).Value = '5.89'
$ws.Cells.Item(25, 5).Value = '  +9.92%  '

$ws.Cells.Item(26, 4).Value = '93.34'
$ws.Cells.Item(26, 5).Value = '  +11.40%  '

$ws.Cells.Item(27, 4).Value = '12.14'
$ws.Cells.Item(27, 5).Value = '  +4.87%  '

$ws.Cells.Item(28, 4).Value = '3.315.63'
$ws.Cells.Item(28, 5).Value = '  +3.56%  '

$ws.Cells.Item(29, 5).Value = '  +0.08%  '

$ws.Cells.Item(30, 5).Value = '  +2.67%  '

$ws.Cells.Item(31, 4).Value = '0.162'
$ws.Cells.Item(31, 5).Value = '  +0.11%  '

$ws.Cells.Item(32, 4).Value = '9.16'
$ws.Cells.Item(32, 5).Value = '  +11.65%  '

$ws.Cells.Item(33, 4).Value = '26.96'
$ws.Cells.Item(33, 5).Value = '  +18.16%  '

$ws.Cells.Item(34, 4).Value = '521.43'
$ws.Cells.Item(34, 5).Value = '  +4.04%  '

$ws.Cells.Item(35, 4).Value = '0.181'
$ws.Cells.Item(35, 5).Value = '  +32.30%  '

$ws.Cells.Item(36, 2).Value = 'dogwifhat'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(36, 4).Value = '3.65'
$ws.Cells.Item(36, 5).Value = '  -1.59%  '

$ws.Cells.Item(37, 2).Value = 'PancakeSwap'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(37, 4).Value = '1.93'
$ws.Cells.Item(37, 5).Value = '  +7.43%  '

$ws.Cells.Item(38, 2).Value = 'RenderToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Cells.Item(38, 4).Value = '6.91'
$ws.Cells.Item(38, 5).Value = '  +4.60%  '

$ws.Cells.Item(39, 2).Value = 'Fetch.AI'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(39, 4).Value = '1.31'
$ws.Cells.Item(39, 5).Value = '  +6.34%  '

$ws.Cells.Item(40, 2).Value = 'Kaspa'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(40, 4).Value = '0.142'
$ws.Cells.Item(40, 5).Value = '  +8.24%  '

$ws.Cells.Item(41, 4).Value = '0.0872'
$ws.Cells.Item(41, 5).Value = '  +27.12%  '

$ws.Cells.Item(42, 5).Value = '  -0.03%  '

$ws.Cells.Item(43, 5).Value = '  +0.08%  '

$ws.Cells.Item(44, 4).Value = '0.414'
$ws.Cells.Item(44, 5).Value = '  +14.78%  '

$ws.Cells.Item(45, 4).Value = '1.98'
$ws.Cells.Item(45, 5).Value = '  +8.47%  '

$ws.Cells.Item(46, 5).Value = '  -0.01%  '

$ws.Cells.Item(47, 4).Value = '150.70'
$ws.Cells.Item(47, 5).Value = '  +2.87%  '

$ws.Cells.Item(48, 4).Value = '45.44'
$ws.Cells.Item(48, 5).Value = '  +5.07%  '

$ws.Cells.Item(49, 2).Value = 'ImmutableX'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(49, 4).Value = '1.36'
$ws.Cells.Item(49, 5).Value = '  +12.27%  '

$ws.Cells.Item(50, 2).Value = 'Filecoin'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(50, 4).Value = '4.54'
$ws.Cells.Item(50, 5).Value = '  +9.96%  '

$ws.Cells.Item(51, 5).Value = '  +15.44%  '
